$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the (empty) row 2 gap row ---------------------------------
# Row 2 has no cell content in the target; touching Hidden (on then off)
# makes the sheet materialize an explicit empty <row r="2"/> element
# without leaving any residual attribute behind.
$ws.Rows.Item(2).Hidden = $true
$ws.Rows.Item(2).Hidden = $false

# --- Clear the two stray empty cells at the end of row 127 -------------
$ws.Range("H127:I127").ClearContents()

# --- Helper: write a literal ISO date string into column A -------------
# Plain `.Value = "2024-05-21"` gets auto-recognized by Excel's type
# inference as a real date (and stamps a date-formatted style on the
# cell). The source file stores these as plain text, so force a Text
# number format while assigning, then reset the cell style back to
# Normal (default) - the value stays text but no residual style/format
# is left behind on the cell.
function Set-TextDate($rng, $value) {
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# --- Append the four new ticket rows (128-131) --------------------------
Set-TextDate $ws.Range("A128") "2024-05-21"
$ws.Range("B128").Value = "11:42:31"
$ws.Range("C128").Value = "No coge placa"
$ws.Range("D128").Value = "-"
$ws.Range("E128").Value = "-"
$ws.Range("F128").Value = "-"
$ws.Range("G128").Value = "-"

Set-TextDate $ws.Range("A129") "2024-05-21"
$ws.Range("B129").Value = "11:42:35"
$ws.Range("C129").Value = "Ascensor no sube"
$ws.Range("D129").Value = "-"
$ws.Range("E129").Value = "-"
$ws.Range("F129").Value = "-"
$ws.Range("G129").Value = "-"
$ws.Range("H129").Value = "11:42:37"
$ws.Range("I129").Value = "0:00:02"

Set-TextDate $ws.Range("A130") "2024-05-21"
$ws.Range("B130").Value = "11:43:12"
$ws.Range("C130").Value = "-"
$ws.Range("D130").Value = "AOI (malla)"
$ws.Range("E130").Value = "-"
$ws.Range("F130").Value = "-"
$ws.Range("G130").Value = "-"
$ws.Range("H130").Value = "11:43:14"
$ws.Range("I130").Value = "0:00:02"

Set-TextDate $ws.Range("A131") "2024-05-21"
$ws.Range("B131").Value = "11:43:16"
$ws.Range("C131").Value = "-"
$ws.Range("D131").Value = "Cámara no detecta foams"
$ws.Range("E131").Value = "-"
$ws.Range("F131").Value = "-"
$ws.Range("G131").Value = "-"
$ws.Range("H131").Value = "11:43:18"
$ws.Range("I131").Value = "0:00:02"
